$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.109.54"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "1.918.65"
$ws.Range("E3").Value = "  +2.41%  "

$ws.Range("E4").Value = "  +0.25%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "319.38"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("E6").Value = "  +0.17%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.5071"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -0.51%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4061"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +2.63%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08332"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("E10").Value = "  +1.83%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "42.05"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("E12").Value = "  +1.23%  "

$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").Value = "1.917.14"
$ws.Range("E14").Value = "  +2.53%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.248"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +0.52%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +0.32%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "92.55"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("E18").Value = "  +0.81%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06497"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +1.66%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.46"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  +2.19%  "

$ws.Range("E21").Value = "  +0.18%  "

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.945"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +1.66%  "

$ws.Range("D23").Value = "30.124.45"
$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("D26").Value = "2.137.00"
$ws.Range("E26").Value = "  +2.49%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.85"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  +3.38%  "

$ws.Range("E28").Value = "  +1.07%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.262"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +0.64%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "128.93"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +0.95%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.132"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  +5.44%  "

$ws.Range("E32").Value = "  +1.04%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.946"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -0.06%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.786"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  +1.95%  "

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.02445"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  +0.11%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.312"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  +1.45%  "

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06432"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +0.81%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.233"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +4.44%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2145"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -0.05%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6462"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  +2.18%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.615"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  +0.72%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.46"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  +0.29%  "

$ws.Range("E43").Value = "  +0.61%  "

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.36"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  +2.70%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6050"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +2.20%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.174"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +7.56%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.626"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -0.41%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "122.17"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -0.47%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.207"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -0.02%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.134"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +0.94%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "78.02"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +0.94%  "

